$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------------
# 1) Reposition the four existing pictures / shapes (sizes stay the same).
#    Shape.Left / Shape.Top are expressed in points; the literals below are
#    the point values whose internal (truncating) EMU conversion reproduces
#    the exact target EMU offsets from the diff.
# ---------------------------------------------------------------------------

# Picture 4  : off x="686873"  y="944629"   -> x="-11185" y="921249"
$picTop = $s.Shapes.Item(2)
$picTop.Left = -0.8807086944580078
$picTop.Top  = 72.53929138183594

# 箭头: 下 8 (down arrow) : off x="4655898" y="2715199" -> x="3615663" y="2715199"
$arrow = $s.Shapes.Item(3)
$arrow.Left = 284.6978759765625
$arrow.Top  = 213.7952117919922

# Picture 10 : off x="1492166" y="3162824" -> x="451931" y="3162824"
$picMid = $s.Shapes.Item(4)
$picMid.Left = 35.585121154785156
$picMid.Top  = 249.04127502441406

# Picture 12 : off x="1405883" y="4104835" -> x="365648" y="4104835"
$picBottom = $s.Shapes.Item(5)
$picBottom.Left = 28.791181564331055
$picBottom.Top  = 323.2153625488281

# ---------------------------------------------------------------------------
# 2) Add the new "Original vs Modified" text box.
# ---------------------------------------------------------------------------
$txt = $s.Shapes.AddTextbox(1, 554.302001953125, 84.55842590332031, 386.8529357910156, 435.6249694824219)
$txt.Name = "文本占位符 3"

$tf = $txt.TextFrame
$tf.AutoSize = 2
$tf.MarginLeft = 7.2
$tf.MarginTop = 3.6
$tf.MarginRight = 7.2
$tf.MarginBottom = 3.6

$lines = @(
    "Original:",
    "Avg seq Length: 163.50",
    "Num of training data: 6040",
    "Num of items: 3416",
    "",
    "Modified:",
    "Avg seq Length: 185.46",
    "Num of training data: 121919",
    "Num of items: 3416"
)
$tf.TextRange.Text = [string]::Join("`r", $lines)

# Text colour (#1F2328) and typeface (-apple-system) for every run, bullets
# removed from every paragraph, plus the default body size (28pt) that this
# shape's style would otherwise inherit.
$textColor = 2630431  # RGB(0x1F, 0x23, 0x28)
for ($i = 1; $i -le $lines.Count; $i++) {
    $line = $tf.TextRange.Lines($i, 1)
    $line.ParagraphFormat.Bullet.Visible = 0
    if ($line.Length -gt 0) {
        $line.Font.Size = 28
        $line.Font.Color.RGB = $textColor
        $line.Font.Name = "-apple-system"
    }
}

# Paragraph 2 ("Avg seq Length" + ": 163.50") is split across two runs in
# the source deck - replicate that run break (cosmetic rPr-only split).
$p2 = $tf.TextRange.Lines(2, 1)
$split = $p2.Characters(1, 14)
$split.Font.Size = 28
$split.Font.Color.RGB = $textColor
$split.Font.Name = "-apple-system"
